$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the aa8163b8-... file
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# zh-cn sheet: row 3 is the aa8163b8-... file
$wsZhCn.Range("C3").Value = $handedBack
$wsZhCn.Range("K3").Value = "2016-08-25 04:46:57"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: row 3 is the aa8163b8-... file
$wsDeDe.Range("C3").Value = $handedBack
$wsDeDe.Range("K3").Value = "2016-08-25 04:47:09"
$wsDeDe.Range("P3").Value = ""
